$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Br"
$ws.Range("B3").Value = 30.46543312072754

$ws.Range("A4").Value = "CaCO3"
$ws.Range("B4").Value = 31.25371170043945
